# ---------------------------------------------------------------------------
# Applies the "grammar and tester rework; added enum instead of attribute
# token" edit described by the supplied diff.
#
# Strategy notes (discovered empirically against this COM-interop runtime):
#   * Range.HighlightColorIndex (set directly on a Range) is bugged in this
#     engine -- it always paints paragraph #1, ignoring the receiver. Going
#     through Range.Font.HighlightColorIndex instead dispatches correctly
#     and also properly fragments runs / stamps the paragraph-mark rPr
#     (w:pPr/w:rPr) when applied to a Paragraph.Range.
#   * Range.InsertAfter / InsertBefore reliably mint brand new dedicated
#     <w:r> runs (no unwanted merging with neighbours), which lets us build
#     up the exact run sequence the diff wants, one text segment at a time.
#   * Toggling a boolean character property away from and back to its
#     *original* value (e.g. Italic 1 -> 0 -> 1) splits the run at the
#     sub-range boundary while leaving the final <w:rPr> identical to what
#     it would have been natively (no stray residue).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ===========================================================================
# 1) Paragraph "Abbiamo inserito troppi TOKEN..." -> split out "esere" so it
#    can be wrapped the way the diff shows (proofErr markers aren't part of
#    the exposed object model, so we focus on getting the run boundaries /
#    text right).
# ===========================================================================
$p2 = $d.Paragraphs(2)
$t2 = $p2.Range.Text
$rel = $t2.IndexOf("esere")
$abs = $p2.Range.Start + $rel
$esereRange = $d.Range($abs, $abs + 5)
$esereRange.Bold = $true
$esereRange2 = $d.Range($abs, $abs + 5)
$esereRange2.Bold = $false

# ===========================================================================
# 2) Paragraph "Dove può essere accettato un ID..."
#    a) split the italic "rservedWordRule" run into "r" + "e" + "servedWordRule"
#    b) append the new sentence about ASSIGN(=) at the end of the paragraph
# ===========================================================================
$p3 = $d.Paragraphs(3)
$t3 = $p3.Range.Text
$rel3 = $t3.IndexOf("rservedWordRule")
$abs3 = $p3.Range.Start + $rel3

# "r" (offset 0..1) and "e" (offset 1..2) each get split out of the italic
# run by toggling Italic off then back on -- leaves a clean <w:i/><w:iCs/>
# on every fragment, matching the original run's formatting exactly.
$rSplit1 = $d.Range($abs3, $abs3 + 1)
$rSplit1.Italic = $false
$rSplit1b = $d.Range($abs3, $abs3 + 1)
$rSplit1b.Italic = $true

$rSplit2 = $d.Range($abs3 + 1, $abs3 + 2)
$rSplit2.Italic = $false
$rSplit2b = $d.Range($abs3 + 1, $abs3 + 2)
$rSplit2b.Italic = $true

# Append the new trailing sentence as three distinct runs (plain text +
# "ASSIGN(" + the remainder), mirroring the diff's run layout.
$p3b = $d.Paragraphs(3)
$p3b.Range.InsertAfter("; anche ")
$p3c = $d.Paragraphs(3)
$p3c.Range.InsertAfter("ASSIGN(")
$p3d = $d.Paragraphs(3)
$p3d.Range.InsertAfter("=) può essere utilizzato negli ID.")

# ===========================================================================
# 3) Paragraph "Inserire le lettere greche (almeno)..." -> highlight the
#    whole paragraph (and its paragraph mark) yellow.
# ===========================================================================
$p4 = $d.Paragraphs(4)
$p4.Range.Font.HighlightColorIndex = 7   # wdYellow

# ===========================================================================
# 4) Paragraph "Inserire tutti i caratteri speciali..."
#    a) merge the "uno" gram-checked fragment back into the surrounding text
#    b) highlight the (now merged) original sentence yellow, incl. the
#       paragraph mark
#    c) append the new red-highlighted aside about special characters
#       breaking the name, plus a yellow ":" / " " separator first.
# ===========================================================================
$p5 = $d.Paragraphs(5)
$t5 = $p5.Range.Text
$relStart = $t5.IndexOf("Inserire tutti i caratteri speciali. Vanno messi ")
$relEnd = $t5.IndexOf(" a uno in un ") + (" a uno in un ").Length
$absStart = $p5.Range.Start + $relStart
$absEnd = $p5.Range.Start + $relEnd
$mergeRange = $d.Range($absStart, $absEnd)
$mergeRange.Text = ""
$insertPoint = $d.Range($absStart, $absStart)
$insertPoint.InsertBefore("Inserire tutti i caratteri speciali. Vanno messi uno a uno in un ")

# Highlight the whole (pre-existing) paragraph content + paragraph mark
# yellow before appending the new red-highlighted text.
$p5b = $d.Paragraphs(5)
$p5b.Range.Font.HighlightColorIndex = 7   # wdYellow

function Append-Highlighted($paraIndex, $text, $color) {
    $p = $d.Paragraphs($paraIndex)
    $startPos = $p.Range.End - 1
    $p.Range.InsertAfter($text)
    $endPos = $startPos + $text.Length
    $r = $d.Range($startPos, $endPos)
    $r.Font.HighlightColorIndex = $color
}

Append-Highlighted 5 ":" 7                                                             # wdYellow
Append-Highlighted 5 " " 7                                                             # wdYellow
Append-Highlighted 5 "alcuni " 6                                                       # wdRed
Append-Highlighted 5 "caratteri da errore nel name (es. parentesi ma invece nel " 6     # wdRed
Append-Highlighted 5 "value" 6                                                         # wdRed
Append-Highlighted 5 " si può. Direi di non farci problemi per ora)" 6                 # wdRed

# ===========================================================================
# 5) Paragraph "Per fare riconoscere ID separati..." -> paragraph mark
#    highlighted yellow, but the run itself highlighted red.
# ===========================================================================
$p6 = $d.Paragraphs(6)
$p6.Range.Font.HighlightColorIndex = 7    # wdYellow (stamps pPr + run)
$p6b = $d.Paragraphs(6)
$runOnly6 = $d.Range($p6b.Range.Start, $p6b.Range.End - 1)
$runOnly6.Font.HighlightColorIndex = 6    # wdRed (overwrite the run only)

# ===========================================================================
# 6) Paragraph "Inserire EOF per risolvere..." -> both paragraph mark and
#    run highlighted yellow.
# ===========================================================================
$p7 = $d.Paragraphs(7)
$p7.Range.Font.HighlightColorIndex = 7    # wdYellow
